$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Sentinel Defined" (column H) values from "No" to "Yes" for rows 2, 4, 7
$ws.Range("H2").Value = "Yes"
$ws.Range("H4").Value = "Yes"
$ws.Range("H7").Value = "Yes"

# Update the active cell selection to I9
$ws.Range("I9").Select()
